$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: force a Range to hold a literal text value (shared string) even when
# the cell's number format would otherwise cause Excel to coerce the string to
# a number (columns like P/Q use numeric-looking formats but store text).
function Set-TextValue($rng, [string]$val) {
    $origFmt = $rng.NumberFormat
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = $origFmt
}

# ---------------------------------------------------------------------------
# 1. Insert a new row above row 13 (shifts the old row 13 "سرنجات" product
#    row down to 14, the totals row down to 15, and the footer row down to 16)
# ---------------------------------------------------------------------------
$ws.Rows("13:13").Insert()
$ws.Rows("13:13").RowHeight = 24.75

# ---------------------------------------------------------------------------
# 2. Copy the cell formatting for every column of the (now shifted) product
#    row 14 into the new, still-empty row 13, so the new row matches the
#    look of the other product rows exactly.
# ---------------------------------------------------------------------------
foreach ($col in @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")) {
    $ws.Range("$col`14").Copy()
    $ws.Range("$col`13").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 3. Re-create the merged cells for the new product row 13
# ---------------------------------------------------------------------------
$ws.Range("A13:B13").Merge()
$ws.Range("C13:G13").Merge()
$ws.Range("H13:K13").Merge()
$ws.Range("L13:M13").Merge()
$ws.Range("N13:O13").Merge()

# ---------------------------------------------------------------------------
# 4. Fill in the new product row 13 with the "PANADOL ADVANCE 500 MG 48
#    TABLETS" line item.
# ---------------------------------------------------------------------------
$ws.Range("A13").Value = 7
Set-TextValue $ws.Range("C13") "PANADOL ADVANCE 500 MG 48 TABLETS"
Set-TextValue $ws.Range("H13") "1:3"
Set-TextValue $ws.Range("L13") "1"
Set-TextValue $ws.Range("N13") "92.00"
Set-TextValue $ws.Range("P13") "23.0000"
Set-TextValue $ws.Range("Q13") "0:1"

# ---------------------------------------------------------------------------
# 5. Update the totals row (now row 15) and the footer timestamp (now row 16)
# ---------------------------------------------------------------------------
$ws.Range("P15").Value = 194
$ws.Range("A16").Value = "Friday, 11 July, 2025 5:39 PM"

Write-Output "Edit applied"
